$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2000
$ws.Range("J7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("N7").Value = -2224
$ws.Range("H13").Value = 5500
$ws.Range("J13").Value = 5500
$ws.Range("L13").Value = 5500
$ws.Range("N13").Value = -5838
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2382
$ws.Range("H33").Value = 134.2
$ws.Range("I33").Value = 134.2
$ws.Range("K33").Value = 134.2
$ws.Range("M33").Value = 94.80000000000001
$ws.Range("H41").Value = 615.9
$ws.Range("I41").Value = 343
$ws.Range("K41").Value = 343
$ws.Range("M41").Value = 97
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H76").Value = 1676989.6
$ws.Range("I76").Value = 2605413.2
$ws.Range("K76").Value = 2605413.2
$ws.Range("M76").Value = -2605098.2
$ws.Range("H79").Value = 1676989.6
$ws.Range("I79").Value = 2605413.2
$ws.Range("K79").Value = 2605413.2
$ws.Range("M79").Value = -2604321.2
$ws.Range("H113").Value = 22330.666
$ws.Range("I113").Value = 24396.8
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 24396.8
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -21142.8
$ws.Range("N113").Value = -18508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1650.9231
$ws.Range("I45").Value = 1255.8
$ws.Range("J45").Value = 1897.875
$ws.Range("K45").Value = 1255.8
$ws.Range("L45").Value = 1897.875
$ws.Range("M45").Value = -878.8
$ws.Range("N45").Value = -2651.875
$ws.Range("H74").Value = 1188.44
$ws.Range("I74").Value = 1000
$ws.Range("K74").Value = 1000
$ws.Range("M74").Value = -126
$ws.Range("H77").Value = 1188.44
$ws.Range("I77").Value = 1000
$ws.Range("K77").Value = 5000
$ws.Range("M77").Value = -632
$ws.Range("H88").Value = 25966.445
$ws.Range("I88").Value = 2949.5
$ws.Range("K88").Value = 2949.5
$ws.Range("M88").Value = -2543.5
$ws.Range("H91").Value = 25966.445
$ws.Range("I91").Value = 2949.5
$ws.Range("K91").Value = 2949.5
$ws.Range("M91").Value = -1545.5
$ws.Range("H102").Value = 2206.5881
$ws.Range("I102").Value = 2034.1333
$ws.Range("K102").Value = 2034.1333
$ws.Range("M102").Value = -412.1333
$ws.Range("H110").Value = 1753.2609
$ws.Range("I110").Value = 1229.0667
$ws.Range("K110").Value = 1229.0667
$ws.Range("M110").Value = 815.9332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 566.88464
$ws.Range("I94").Value = 532.43475
$ws.Range("K94").Value = 532.43475
$ws.Range("M94").Value = -81.43475000000001
$ws.Range("H108").Value = 94980
$ws.Range("J108").Value = 94980
$ws.Range("L108").Value = 94980
$ws.Range("N108").Value = -102660

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 2504.1667
$ws.Range("I3").Value = 1005
$ws.Range("K3").Value = 1005
$ws.Range("M3").Value = -892
$ws.Range("H31").Value = 2799.3
$ws.Range("I31").Value = 1741.25
$ws.Range("K31").Value = 1741.25
$ws.Range("M31").Value = -1446.25
$ws.Range("H34").Value = 2799.3
$ws.Range("I34").Value = 1741.25
$ws.Range("K34").Value = 1741.25
$ws.Range("M34").Value = -1539.25
$ws.Range("H58").Value = 1554353.9
$ws.Range("I58").Value = 3953827.5
$ws.Range("J58").Value = 1753.2941
$ws.Range("K58").Value = 3953827.5
$ws.Range("L58").Value = 1753.2941
$ws.Range("M58").Value = -3953624.5
$ws.Range("N58").Value = -2159.2941
$ws.Range("H99").Value = 418738.28
$ws.Range("I99").Value = 715707.9
$ws.Range("J99").Value = 2980.9
$ws.Range("K99").Value = 715707.9
$ws.Range("L99").Value = 2980.9
$ws.Range("M99").Value = -714209.9
$ws.Range("N99").Value = -5976.9
$ws.Range("H126").Value = 418738.28
$ws.Range("I126").Value = 715707.9
$ws.Range("J126").Value = 2980.9
$ws.Range("K126").Value = 2147123.7
$ws.Range("L126").Value = 8942.700000000001
$ws.Range("M126").Value = -2144653.7
$ws.Range("N126").Value = -13882.7
$ws.Range("H133").Value = 61932.332
$ws.Range("J133").Value = 61932.332
$ws.Range("L133").Value = 61932.332
$ws.Range("N133").Value = -66992.33199999999
$ws.Range("H136").Value = 1554353.9
$ws.Range("I136").Value = 3953827.5
$ws.Range("J136").Value = 1753.2941
$ws.Range("K136").Value = 11861482.5
$ws.Range("L136").Value = 5259.8823
$ws.Range("M136").Value = -11858932.5
$ws.Range("N136").Value = -10359.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 4109.4443
$ws.Range("I109").Value = 999.7143
$ws.Range("J109").Value = 6088.364
$ws.Range("K109").Value = 2999.1429
$ws.Range("L109").Value = 18265.092
$ws.Range("M109").Value = -1959.1429
$ws.Range("N109").Value = -20345.092
$ws.Range("H131").Value = 19726.666
$ws.Range("J131").Value = 22908.195
$ws.Range("L131").Value = 68724.58499999999
$ws.Range("N131").Value = -78804.58499999999
$ws.Range("H134").Value = 2498.9443
$ws.Range("I134").Value = 1922.4615
$ws.Range("K134").Value = 5767.3845
$ws.Range("M134").Value = -697.3845000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 1569.5
$ws.Range("I97").Value = 1263.3334
$ws.Range("K97").Value = 1263.3334
$ws.Range("M97").Value = -767.3334
$ws.Range("H113").Value = 1449.75
$ws.Range("J113").Value = 1599.6666
$ws.Range("L113").Value = 1599.6666
$ws.Range("N113").Value = -5939.6666
$ws.Range("H126").Value = 2830310.2
$ws.Range("I126").Value = 3270424
$ws.Range("K126").Value = 9811272
$ws.Range("M126").Value = -9808802
$ws.Range("H132").Value = 4812356
$ws.Range("I132").Value = 6412360
$ws.Range("J132").Value = 12344.5
$ws.Range("K132").Value = 19237080
$ws.Range("L132").Value = 37033.5
$ws.Range("M132").Value = -19234550
$ws.Range("N132").Value = -42093.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2972.647
$ws.Range("I16").Value = 3364.9
$ws.Range("J16").Value = 2412.2856
$ws.Range("K16").Value = 3364.9
$ws.Range("L16").Value = 2412.2856
$ws.Range("M16").Value = -3194.9
$ws.Range("N16").Value = -2752.2856
$ws.Range("H40").Value = 9800
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H46").Value = 2092.3076
$ws.Range("J46").Value = 2327.182
$ws.Range("L46").Value = 2327.182
$ws.Range("N46").Value = -2703.182
$ws.Range("H61").Value = 4614.2856
$ws.Range("I61").Value = 5125
$ws.Range("J61").Value = 3933.3333
$ws.Range("K61").Value = 5125
$ws.Range("L61").Value = 3933.3333
$ws.Range("M61").Value = -4923
$ws.Range("N61").Value = -4337.3333
$ws.Range("H113").Value = 4614.2856
$ws.Range("I113").Value = 5125
$ws.Range("J113").Value = 3933.3333
$ws.Range("K113").Value = 5125
$ws.Range("L113").Value = 3933.3333
$ws.Range("M113").Value = -2955
$ws.Range("N113").Value = -8273.3333
$ws.Range("H122").Value = 10335.5
$ws.Range("I122").Value = 7336
$ws.Range("K122").Value = 22008
$ws.Range("M122").Value = -19558

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H62").Value = 2849.5
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 2899
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 2899
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -4147
$ws.Range("H65").Value = 2849.5
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 2899
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 14495
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -20735
$ws.Range("H122").Value = 100277.52
$ws.Range("I122").Value = 118628.625
$ws.Range("K122").Value = 355885.875
$ws.Range("M122").Value = -353435.875
$ws.Range("H126").Value = 8043.2383
$ws.Range("I126").Value = 7579.5835
$ws.Range("J126").Value = 8661.444
$ws.Range("K126").Value = 22738.7505
$ws.Range("L126").Value = 25984.332
$ws.Range("M126").Value = -20268.7505
$ws.Range("N126").Value = -30924.332
$ws.Range("H132").Value = 1859.6
$ws.Range("I132").Value = 1386
$ws.Range("K132").Value = 4158
$ws.Range("M132").Value = -1628
